$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 29600
$ws.Range("J18").Value = 21500
$ws.Range("L18").Value = 21500
$ws.Range("N18").Value = -22068

$ws.Range("H31").Value = 3
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 9
$ws.Range("M31").Value = 221

$ws.Range("H32").Value = 10250
$ws.Range("I32").Value = 8000
$ws.Range("K32").Value = 8000
$ws.Range("M32").Value = -7674

$ws.Range("H40").Value = 4218.625
$ws.Range("I40").Value = 2374
$ws.Range("J40").Value = 4482.143
$ws.Range("K40").Value = 2374
$ws.Range("L40").Value = 4482.143
$ws.Range("M40").Value = -2199
$ws.Range("N40").Value = -4832.143

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H53").Value = 271.13333
$ws.Range("I53").Value = 264.33334
$ws.Range("K53").Value = 264.33334
$ws.Range("M53").Value = 372.66666

$ws.Range("H55").Value = 291.66666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 291.66666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 291.66666
$ws.Range("N55").Value = -719.66666
$ws.Range("M55").ClearContents()

$ws.Range("I69").Value = 7015
$ws.Range("K69").Value = 21045
$ws.Range("M69").Value = -20171

$ws.Range("H70").Value = 18966.5
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 27749.75
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 83249.25
$ws.Range("M70").Value = -3930
$ws.Range("N70").Value = -83789.25

$ws.Range("I72").Value = 7015
$ws.Range("K72").Value = 63135
$ws.Range("M72").Value = -58767

$ws.Range("H73").Value = 18966.5
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 27749.75
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 83249.25
$ws.Range("M73").Value = -3264
$ws.Range("N73").Value = -85121.25

$ws.Range("H113").Value = 3166.3333
$ws.Range("I113").Value = 2999.5
$ws.Range("K113").Value = 2999.5
$ws.Range("M113").Value = 254.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1008
$ws.Range("I4").Value = 812
$ws.Range("K4").Value = 812
$ws.Range("M4").Value = -696

$ws.Range("H5").Value = 18.875
$ws.Range("I5").Value = 10.333333
$ws.Range("J5").Value = 24
$ws.Range("K5").Value = 10.333333
$ws.Range("L5").Value = 24
$ws.Range("M5").Value = 101.666667
$ws.Range("N5").Value = -248

$ws.Range("H6").Value = 10000500
$ws.Range("I6").Value = 10000500
$ws.Range("K6").Value = 10000500
$ws.Range("M6").Value = -10000327

$ws.Range("H95").Value = 13500
$ws.Range("J95").Value = 13500
$ws.Range("L95").Value = 13500
$ws.Range("N95").Value = -18992

$ws.Range("H96").Value = 14000
$ws.Range("J96").Value = 14000
$ws.Range("L96").Value = 14000
$ws.Range("N96").Value = -19492

$ws.Range("H122").Value = 2384.7144
$ws.Range("I122").Value = 2148.25
$ws.Range("K122").Value = 6444.75
$ws.Range("M122").Value = -3994.75

$ws.Range("H132").Value = 2981.8572
$ws.Range("I132").Value = 2714.75
$ws.Range("K132").Value = 8144.25
$ws.Range("M132").Value = -5614.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 18.875
$ws.Range("I4").Value = 10.333333
$ws.Range("J4").Value = 24
$ws.Range("K4").Value = 10.333333
$ws.Range("L4").Value = 24
$ws.Range("M4").Value = 104.666667
$ws.Range("N4").Value = -254

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H22").Value = 813.7143
$ws.Range("I22").Value = 782.6667
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 782.6667
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -609.6667
$ws.Range("N22").Value = -1346

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H86").Value = 1299.6666
$ws.Range("I86").Value = 1249.5
$ws.Range("J86").Value = 1400
$ws.Range("K86").Value = 1249.5
$ws.Range("L86").Value = 1400
$ws.Range("M86").Value = -126.5
$ws.Range("N86").Value = -3646

$ws.Range("H89").Value = 1299.6666
$ws.Range("I89").Value = 1249.5
$ws.Range("J89").Value = 1400
$ws.Range("K89").Value = 6247.5
$ws.Range("L89").Value = 7000
$ws.Range("M89").Value = -631.5
$ws.Range("N89").Value = -18232

$ws.Range("H103").Value = 3477.75
$ws.Range("J103").Value = 3477.75
$ws.Range("L103").Value = 3477.75
$ws.Range("N103").Value = -5821.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 162.25
$ws.Range("I19").Value = 162.25
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 162.25
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 7.75
$ws.Range("N19").ClearContents()

$ws.Range("H22").Value = 281
$ws.Range("I22").Value = 281
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 281
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 69
$ws.Range("N22").ClearContents()

$ws.Range("H24").Value = 162.25
$ws.Range("I24").Value = 162.25
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 162.25
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 7.75
$ws.Range("N24").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H132").Value = 2599.75
$ws.Range("I132").Value = 2599.75
$ws.Range("K132").Value = 7799.25
$ws.Range("M132").Value = -5269.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1903.2
$ws.Range("I98").Value = 1887.5
$ws.Range("J98").Value = 1913.6666
$ws.Range("K98").Value = 5662.5
$ws.Range("L98").Value = 5740.9998
$ws.Range("M98").Value = -4164.5
$ws.Range("N98").Value = -8736.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5742
$ws.Range("I126").Value = 4470
$ws.Range("K126").Value = 13410
$ws.Range("M126").Value = -10940

$ws.Range("H132").Value = 8798.6
$ws.Range("I132").Value = 7998
$ws.Range("K132").Value = 23994
$ws.Range("M132").Value = -21464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9099.6
$ws.Range("I7").Value = 8749.5
$ws.Range("J7").Value = 9333
$ws.Range("K7").Value = 8749.5
$ws.Range("L7").Value = 9333
$ws.Range("M7").Value = -8637.5
$ws.Range("N7").Value = -9557

$ws.Range("H46").Value = 4730.6665
$ws.Range("I46").Value = 3262.6667
$ws.Range("J46").Value = 4864.121
$ws.Range("K46").Value = 3262.6667
$ws.Range("L46").Value = 4864.121
$ws.Range("M46").Value = -3074.6667
$ws.Range("N46").Value = -5240.121

$ws.Range("H56").Value = 20050
$ws.Range("I56").Value = 20050
$ws.Range("K56").Value = 20050
$ws.Range("M56").Value = -19359

$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

$ws.Range("H82").Value = 1979.9231
$ws.Range("J82").Value = 1848.625
$ws.Range("L82").Value = 1848.625
$ws.Range("N82").Value = -2570.625

$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

$ws.Range("H85").Value = 1979.9231
$ws.Range("J85").Value = 1848.625
$ws.Range("L85").Value = 1848.625
$ws.Range("N85").Value = -4344.625

$ws.Range("H100").Value = 3873.1428
$ws.Range("I100").Value = 4185.3335
$ws.Range("K100").Value = 4185.3335
$ws.Range("M100").Value = -3644.3335

$ws.Range("H126").Value = 9099.6
$ws.Range("I126").Value = 8749.5
$ws.Range("J126").Value = 9333
$ws.Range("K126").Value = 26248.5
$ws.Range("L126").Value = 27999
$ws.Range("M126").Value = -23778.5
$ws.Range("N126").Value = -32939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H81").Value = 37750.75
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 37750.75
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H97").Value = 50000
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982

$ws.Range("H101").Value = 23044
$ws.Range("J101").Value = 23044
$ws.Range("L101").Value = 23044
$ws.Range("N101").Value = -29534

$ws.Range("H122").Value = 2760.4285
$ws.Range("I122").Value = 2760.4285
$ws.Range("K122").Value = 8281.2855
$ws.Range("M122").Value = -5831.2855

$ws.Range("H126").Value = 2016
$ws.Range("I126").Value = 2016
$ws.Range("K126").Value = 6048
$ws.Range("M126").Value = -3578
